$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# --- Row 151: BCR register header, field BS ---
$ws.Range("A151").Value = "BCR"
$ws.Range("B151").Value = 0
$ws.Range("C151").Value = "BS"
$ws.Range("D151").Value = 0
$ws.Range("E151").Value = 16
$ws.Range("F151").Value = 0
$ws.Range("H151").Formula = '="class IOPDmacChannelRegister_"&A151&"_t;"'
$ws.Range("I151").Formula = '="static constexpr u8 "&C151&" = "&B151&";"'
$ws.Range("J151").Formula = '="registerField(Fields::"&C151&", """&C151&""", "&D151&", "&E151&", "&F151&");"'

# --- Row 152: BCR field BA ---
$ws.Range("B152").Value = 1
$ws.Range("C152").Value = "BA"
$ws.Range("D152").Value = 16
$ws.Range("E152").Value = 16
$ws.Range("F152").Value = 0
$ws.Range("I152").Formula = '="static constexpr u8 "&C152&" = "&B152&";"'
$ws.Range("J152").Formula = '="registerField(Fields::"&C152&", """&C152&""", "&D152&", "&E152&", "&F152&");"'

# --- Row 154: CHCR register header, field DR ---
$ws.Range("A154").Value = "CHCR"
$ws.Range("B154").Value = 0
$ws.Range("C154").Value = "DR"
$ws.Range("D154").Value = 0
$ws.Range("E154").Value = 1
$ws.Range("F154").Value = 0
$ws.Range("H154").Formula = '="class IOPDmacChannelRegister_"&A154&"_t;"'
$ws.Range("I154").Formula = '="static constexpr u8 "&C154&" = "&B154&";"'
$ws.Range("J154").Formula = '="registerField(Fields::"&C154&", """&C154&""", "&D154&", "&E154&", "&F154&");"'

# --- Row 155: CHCR field CO ---
$ws.Range("B155").Value = 1
$ws.Range("C155").Value = "CO"
$ws.Range("D155").Value = 9
$ws.Range("E155").Value = 1
$ws.Range("F155").Value = 0
$ws.Range("I155").Formula = '="static constexpr u8 "&C155&" = "&B155&";"'
$ws.Range("J155").Formula = '="registerField(Fields::"&C155&", """&C155&""", "&D155&", "&E155&", "&F155&");"'

# --- Row 156: CHCR field LI ---
$ws.Range("B156").Value = 2
$ws.Range("C156").Value = "LI"
$ws.Range("D156").Value = 10
$ws.Range("E156").Value = 1
$ws.Range("F156").Value = 0
$ws.Range("I156").Formula = '="static constexpr u8 "&C156&" = "&B156&";"'
$ws.Range("J156").Formula = '="registerField(Fields::"&C156&", """&C156&""", "&D156&", "&E156&", "&F156&");"'

# --- Row 157: CHCR field TR ---
$ws.Range("B157").Value = 3
$ws.Range("C157").Value = "TR"
$ws.Range("D157").Value = 24
$ws.Range("E157").Value = 1
$ws.Range("F157").Value = 0
$ws.Range("I157").Formula = '="static constexpr u8 "&C157&" = "&B157&";"'
$ws.Range("J157").Formula = '="registerField(Fields::"&C157&", """&C157&""", "&D157&", "&E157&", "&F157&");"'

# --- Match final selection state from the authored workbook ---
$ws.Range("J154:J157").Select()
